$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 0.09

$ws.Range("B3").Value = 0.07627777677142858
$ws.Range("C3").Value = 0.3935134577428572
$ws.Range("D3").Value = 0.35
$ws.Range("E3").Value = -0.3172356809714286

$ws.Range("B4").Value = 0.07119323431914894
$ws.Range("C4").Value = 0.3856232987021276
$ws.Range("D4").Value = 0.47
$ws.Range("E4").Value = -0.3144300643829787

$ws.Range("D6").Value = 0.09

$ws.Range("B7").Value = 0.05555494788542022
$ws.Range("C7").Value = 0.3858197931862442
$ws.Range("D7").Value = 0.5690000000000001
$ws.Range("E7").Value = -0.330264845300824
$ws.Range("F7").Value = 0.04107092963029224

$ws.Range("B8").Value = 0.06624441533390257
$ws.Range("C8").Value = 0.3656339959397661
$ws.Range("D8").Value = 1.0545
$ws.Range("E8").Value = -0.2993895806058635
$ws.Range("F8").Value = -0.04783411473909105

$ws.Range("B9").Value = 0.07364051646380226
$ws.Range("C9").Value = 0.361513997363836
$ws.Range("D9").Value = 0.835
$ws.Range("E9").Value = -0.2878734809000338
$ws.Range("F9").Value = 0.2313605582209974

$ws.Range("B10").Value = 0.08117092497446987
$ws.Range("C10").Value = 0.3495041029590277
$ws.Range("D10").Value = 0.794
$ws.Range("E10").Value = -0.2683331779845579
$ws.Range("F10").Value = -0.1048885798476854

$ws.Range("B11").Value = 0.07550426201193439
$ws.Range("C11").Value = 0.3645063282456978
$ws.Range("D11").Value = 2.7225
$ws.Range("E11").Value = -0.2890020662337635
$ws.Range("F11").Value = -0.1249384536506635

$ws.Range("B12").Value = 0.07346585858598161
$ws.Range("C12").Value = 0.3514840429189443
$ws.Range("D12").Value = 2.9125
$ws.Range("E12").Value = -0.2780181843329627
$ws.Range("F12").Value = -0.07138323327636287

$ws.Range("B13").Value = 0.07629891094722693
$ws.Range("C13").Value = 0.3736113186116612
$ws.Range("D13").Value = 3.1485
$ws.Range("E13").Value = -0.2973124076644343
$ws.Range("F13").Value = 0.03278845531338925

$ws.Range("B14").Value = 0.07562478363606376
$ws.Range("C14").Value = 0.3715840027926106
$ws.Range("D14").Value = 3.124
$ws.Range("E14").Value = -0.2959592191565468
$ws.Range("F14").Value = 0.1029542503073502

$ws.Range("B15").Value = 0.07069165368653375
$ws.Range("C15").Value = 0.3486501773125208
$ws.Range("D15").Value = 3.1785
$ws.Range("E15").Value = -0.2779585236259871
$ws.Range("F15").Value = -0.03821267699464703

$ws.Range("B16").Value = 0.06754177719287266
$ws.Range("C16").Value = 0.3401305096223162
$ws.Range("D16").Value = 3.002
$ws.Range("E16").Value = -0.2725887324294435
$ws.Range("F16").Value = -0.01952912510577642

$ws.Range("B17").Value = 0.0751396824756658
$ws.Range("C17").Value = 0.345277971315771
$ws.Range("D17").Value = 3.192
$ws.Range("E17").Value = -0.2701382888401052
$ws.Range("F17").Value = -0.09139920879117702

$ws.Range("B18").Value = 0.07729095589413372
$ws.Range("C18").Value = 0.3458161253310973
$ws.Range("D18").Value = 3.0155
$ws.Range("E18").Value = -0.2685251694369636
$ws.Range("F18").Value = -0.09269537133449468

$ws.Range("B19").Value = 0.08603266388702674
$ws.Range("C19").Value = 0.3358995967734651
$ws.Range("D19").Value = 3.242
$ws.Range("E19").Value = -0.2498669328864384
$ws.Range("F19").Value = -0.1010639658503435

$ws.Range("B20").Value = 0.08913090583996927
$ws.Range("C20").Value = 0.3417798123859237
$ws.Range("D20").Value = 3.0885
$ws.Range("E20").Value = -0.2526489065459544
$ws.Range("F20").Value = -0.07314985364866566

$ws.Range("B21").Value = 0.08050259506298765
$ws.Range("C21").Value = 0.3336545065106141
$ws.Range("D21").Value = 2.084
$ws.Range("E21").Value = -0.2531519114476265
$ws.Range("F21").Value = -0.06288030277164036

$ws.Range("B22").Value = 0.08197292239716884
$ws.Range("C22").Value = 0.3156616963791157
$ws.Range("D22").Value = 1.2365
$ws.Range("E22").Value = -0.2336887739819469
$ws.Range("F22").Value = -0.1297323283625916
